$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New handoff batch: two additional localized files finished "Ready for
# handoff" at 2016-03-13 14:40 (zh-cn) / 14:40 (de-de). Append one row per
# file to each of the three sheets, mirroring the existing rows' layout,
# values, and hyperlinks.
# ---------------------------------------------------------------------------

$files = @(
    @{
        Name       = "c86cf631-726d-415a-9e5e-3deb6a4488fa.md"
        HandoffSha = "5f7b1f98e944d6e42f516e17e2850c6262fb5cd1"
        ZhDate     = "2016-03-13 14:40:28"
        DeDate     = "2016-03-13 14:40:31"
        OverviewDate = "2016-40-13 14:40:31"
    },
    @{
        Name       = "efc0dffd-8b80-4f8d-9537-7294ad241cab.md"
        HandoffSha = "ea94803140aa61ff2e0f8ea2829f889980d0325b"
        ZhDate     = "2016-03-13 14:40:28"
        DeDate     = "2016-03-13 14:40:31"
        OverviewDate = "2016-40-13 14:40:31"
    }
)

# --- Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date ------
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewRow = 6
foreach ($f in $files) {
    $r = $overviewRow
    $wsOverview.Range("A$r").Value = $f.Name
    $wsOverview.Hyperlinks.Add($wsOverview.Range("A$r"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$($f.Name)", "", "", $f.Name)
    $wsOverview.Range("B$r").Value = "Ready for handoff"
    $wsOverview.Range("C$r").Value = "Ready for handoff"
    $wsOverview.Range("D$r").Value = $f.OverviewDate
    $overviewRow = $overviewRow + 1
}

# --- Sheets "zh-cn" / "de-de": detail columns -------------------------------
$langSheets = @(
    @{ SheetName = "zh-cn"; Locale = "zh-cn"; DateField = "ZhDate" },
    @{ SheetName = "de-de"; Locale = "de-de"; DateField = "DeDate" }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.SheetName)
    $row = 6
    foreach ($f in $files) {
        $r = $row
        $baseName = $f.Name.Substring(0, $f.Name.Length - 3)
        $handoffFile = "$baseName.$($f.HandoffSha).$($lang.Locale).xlf"
        $handoffDate = $f[$lang.DateField]

        $ws.Range("A$r").Value = $f.Name
        $ws.Hyperlinks.Add($ws.Range("A$r"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$($f.Name)", "", "", $f.Name)

        $ws.Range("B$r").Value = ".md"
        $ws.Hyperlinks.Add($ws.Range("B$r"), "https://github.com/OpenLocalizationTestOrg/oltest.$($lang.Locale)/blob/master/e2e/$($f.Name)", "", "", ".md")

        $ws.Range("C$r").Value = "Ready for handoff"

        $ws.Range("D$r").Value = $handoffFile
        $ws.Hyperlinks.Add($ws.Range("D$r"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.$($lang.Locale)/ci/ht/$handoffFile", "", "", $handoffFile)

        $ws.Range("E$r").Value = $handoffDate
        $ws.Range("H$r").Value = "0001-01-01 00:00:00"
        $ws.Range("I$r").Value = "Include"

        $row = $row + 1
    }
}
